# Updated with pressure and flow sensors
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be created in the same order they appear in the
# target workbook: B10, B11, N10, N11.
$ws.Range("B10").Value2 = "pressure sensor"
$ws.Range("B11").Value2 = "flow sensor"
$ws.Range("N10").Value2 = "car fuel tank pressure sensor"
$ws.Range("N11").Value2 = "car mass air flow sensor"

# --- Row 10: pressure sensor (hardware column B/C, car-parts column N/O/P) ---
$ws.Range("C10").Value2 = 2
$ws.Range("O10").Value2 = 45
$ws.Range("P10").Formula = "=O10*C10"

# --- Row 11: flow sensor (hardware column B/C, car-parts column N/O/P) ---
$ws.Range("C11").Value2 = 2
$ws.Range("O11").Value2 = 35
$ws.Range("P11").Formula = "=O11*C11"

# --- View/selection state: scroll the frozen pane toward the new rows ---
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 10
$ws.Range("O11").Select()
